$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 416, shifting existing rows 416:510 down to 417:511
$ws.Rows(416).Insert()

# Populate the new row 416 with the new price record
$ws.Range("A416").Value = 4
$ws.Range("B416").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C416").Value = "Los Lagos"
$ws.Range("D416").Value = 45173
$ws.Range("E416").Value = 10
$ws.Range("F416").Value = 100112037
$ws.Range("G416").Value = "Cebollín"
$ws.Range("H416").Value = "Sin especificar"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 40
$ws.Range("K416").Value = 6500
$ws.Range("L416").Value = 6500
$ws.Range("M416").Value = 6500
$ws.Range("N416").Value = "$/paquete 36 unidades"
$ws.Range("O416").Value = "Región Metropolitana"
$ws.Range("P416").Value = 181
$ws.Range("Q416").Value = 36
$ws.Range("R416").Value = "Hortaliza"
